$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "AUC" column (E), cloning formatting from column D (F1/ACC-style
#     numeric columns) so the new cells share the existing style entries instead of
#     creating duplicates. ---
$ws.Range("D1:D9").Copy() | Out-Null
$ws.Range("E1:E9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Header row ---
$ws.Range("E1").Value = "AUC"

# --- Data rows: refreshed ACC / Recall / F1 values plus the new AUC column ---
$ws.Range("B2").Value = 0.9070208728652751
$ws.Range("C2").Value = 0.9070208728652751
$ws.Range("D2").Value = 0.9064425894951549
$ws.Range("E2").Value = 0.993093831198981

$ws.Range("B3").Value = 0.5196078431372549
$ws.Range("C3").Value = 0.5196078431372549
$ws.Range("D3").Value = 0.3851106258813493
$ws.Range("E3").Value = 0.8676916925004734

$ws.Range("B4").Value = 0.857685009487666
$ws.Range("C4").Value = 0.857685009487666
$ws.Range("D4").Value = 0.8643479903394946
$ws.Range("E4").Value = 0.9627045973918885

$ws.Range("B5").Value = 0.9120809614168248
$ws.Range("C5").Value = 0.9120809614168248
$ws.Range("D5").Value = 0.911162618569057
$ws.Range("E5").Value = 0.9953297212902542

$ws.Range("B6").Value = 0.9117647058823529
$ws.Range("C6").Value = 0.9117647058823529
$ws.Range("D6").Value = 0.9108786583476612
$ws.Range("E6").Value = 0.9943102371508604

$ws.Range("B7").Value = 0.8956356736242884
$ws.Range("C7").Value = 0.8956356736242884
$ws.Range("D7").Value = 0.8961795097186297
$ws.Range("E7").Value = 0.9871195519191381

$ws.Range("B8").Value = 0.8526249209361164
$ws.Range("C8").Value = 0.8526249209361164
$ws.Range("D8").Value = 0.8518567502077986
$ws.Range("E8").Value = 0.9841194560344118

$ws.Range("B9").Value = 0.3716002530044276
$ws.Range("C9").Value = 0.3716002530044276
$ws.Range("D9").Value = 0.3399893085260272
$ws.Range("E9").Value = 0.6332914087797556

# --- Column widths: column A now matches the other (B/C/D/E) columns' width, and
#     the new column E gets the same width too. ---
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(4).ColumnWidth
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth
